$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 48.8
$ws.Range("J4").Value = 69
$ws.Range("L4").Value = 69
$ws.Range("N4").Value = -297
$ws.Range("H12").Value = 162.5
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -540
$ws.Range("H15").Value = 336.7647
$ws.Range("I15").Value = 336.7647
$ws.Range("K15").Value = 1010.2941
$ws.Range("M15").Value = -841.2941000000001
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -532
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -766
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 876
$ws.Range("I29").Value = 535
$ws.Range("K29").Value = 1605
$ws.Range("M29").Value = -1324
$ws.Range("H32").Value = 1900
$ws.Range("I32").Value = 1900
$ws.Range("K32").Value = 1900
$ws.Range("M32").Value = -1574
$ws.Range("H41").Value = 486
$ws.Range("I41").Value = 377.8
$ws.Range("J41").Value = 666.3333
$ws.Range("K41").Value = 377.8
$ws.Range("L41").Value = 666.3333
$ws.Range("M41").Value = 62.19999999999999
$ws.Range("N41").Value = -1546.3333
$ws.Range("H80").Value = 898.5
$ws.Range("J80").Value = 898.5
$ws.Range("L80").Value = 2695.5
$ws.Range("N80").Value = -4691.5
$ws.Range("H83").Value = 898.5
$ws.Range("J83").Value = 898.5
$ws.Range("L83").Value = 8086.5
$ws.Range("N83").Value = -18070.5
$ws.Range("H116").Value = 9996.25
$ws.Range("J116").Value = 9996.25
$ws.Range("L116").Value = 9996.25
$ws.Range("N116").Value = -16880.25
$ws.Range("H121").Value = 3276.25
$ws.Range("J121").Value = 3276.25
$ws.Range("L121").Value = 9828.75
$ws.Range("N121").Value = -13322.75
$ws.Range("H127").Value = 1137.6666
$ws.Range("J127").Value = 1206.5
$ws.Range("L127").Value = 3619.5
$ws.Range("N127").Value = -13539.5
$ws.Range("H135").Value = 818.8570999999999
$ws.Range("I135").Value = 546.4
$ws.Range("K135").Value = 4917.599999999999
$ws.Range("M135").Value = -2382.599999999999
$ws.Range("H137").Value = 841625.2
$ws.Range("I137").Value = 2006000.4
$ws.Range("J137").Value = 9928.571
$ws.Range("K137").Value = 6018001.199999999
$ws.Range("L137").Value = 29785.713
$ws.Range("M137").Value = -6015451.199999999
$ws.Range("N137").Value = -34885.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 779.5
$ws.Range("I30").Value = 779.5
$ws.Range("K30").Value = 779.5
$ws.Range("M30").Value = -629.5
$ws.Range("H36").Value = 5506.3335
$ws.Range("I36").Value = 5250
$ws.Range("J36").Value = 6019
$ws.Range("K36").Value = 5250
$ws.Range("L36").Value = 6019
$ws.Range("M36").Value = -4904
$ws.Range("N36").Value = -6711
$ws.Range("H92").Value = 95000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H135").Value = 29000
$ws.Range("J135").Value = 29000
$ws.Range("L135").Value = 29000
$ws.Range("N135").Value = -39140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 447.83334
$ws.Range("I22").Value = 266.33334
$ws.Range("J22").Value = 629.3333
$ws.Range("K22").Value = 266.33334
$ws.Range("L22").Value = 629.3333
$ws.Range("M22").Value = 83.66665999999998
$ws.Range("N22").Value = -1329.3333
$ws.Range("H50").Value = 3000
$ws.Range("I50").Value = 3000
$ws.Range("K50").Value = 3000
$ws.Range("M50").Value = -2375
$ws.Range("H51").Value = 30000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 30000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -31472
$ws.Range("H61").Value = 30000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -30696
$ws.Range("H134").Value = 10419.6
$ws.Range("J134").Value = 16166.667
$ws.Range("L134").Value = 48500.001
$ws.Range("N134").Value = -53570.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 95.57143000000001
$ws.Range("I12").Value = 45.666668
$ws.Range("K12").Value = 137.000004
$ws.Range("M12").Value = 35.99999600000001
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.25
$ws.Range("I2").Value = 58.444443
$ws.Range("K2").Value = 58.444443
$ws.Range("M2").Value = 54.555557
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 17504.5
$ws.Range("I36").Value = 17504.5
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 17504.5
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -17019.5
$ws.Range("N36").ClearContents()
$ws.Range("H107").Value = 1499.25
$ws.Range("I107").Value = 1001
$ws.Range("J107").Value = 1997.5
$ws.Range("K107").Value = 1001
$ws.Range("L107").Value = 1997.5
$ws.Range("M107").Value = 919
$ws.Range("N107").Value = -5837.5
$ws.Range("H132").Value = 125992
$ws.Range("I132").Value = 139857.6
$ws.Range("K132").Value = 419572.8
$ws.Range("M132").Value = -417042.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240
$ws.Range("H122").Value = 2362
$ws.Range("I122").Value = 1725
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 5175
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -2725
$ws.Range("N122").Value = -13897
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 9071.286
$ws.Range("I132").Value = 8124.75
$ws.Range("J132").Value = 10333.333
$ws.Range("K132").Value = 24374.25
$ws.Range("L132").Value = 30999.999
$ws.Range("M132").Value = -21844.25
$ws.Range("N132").Value = -36059.999
$ws.Range("H141").Value = 99997.664
$ws.Range("J141").Value = 99997.664
$ws.Range("L141").Value = 99997.664
$ws.Range("N141").Value = -110357.664
